$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update OUTING_TYPE values list (B32:B34) to add low-energy / high-energy
$ws.Range("B32").Value = "fun, exciting, social, chill, relaxing, laid-back, low-energy, high-energy"
$ws.Range("B33").Value = "fun, exciting, social, chill, relaxing, laid-back, low-energy, high-energy"
$ws.Range("B34").Value = "fun, exciting, social, chill, relaxing, laid-back, low-energy, high-energy"

# 2. Update BUDGET,OUTING_TYPE combined entity values (B26:B28) to add low-energy / high-energy combos
$ws.Range("B26").Value = "(low fun),(low exciting),(low social),(low chill),(low relaxing),(low laid-back), (low low-energy), (low high-energy), (high fun),(high exciting),(high social),(high chill),(high relaxing),(high laid-back), (high low-energy), (high high-energy)"
$ws.Range("B27").Value = "(low fun),(low exciting),(low social),(low chill),(low relaxing),(low laid-back), (low low-energy), (low high-energy), (high fun),(high exciting),(high social),(high chill),(high relaxing),(high laid-back), (high low-energy), (high high-energy)"
$ws.Range("B28").Value = "(low fun),(low exciting),(low social),(low chill),(low relaxing),(low laid-back), (low low-energy), (low high-energy), (high fun),(high exciting),(high social),(high chill),(high relaxing),(high laid-back), (high low-energy), (high high-energy)"

# 3. Add new rows for an ADDRESS example and a PHONE_NUMBER entity block
$ws.Range("A35").Value = "`"My address is 123 Westwind.`""
$ws.Range("B35").Clear()

$ws.Range("A36").Value = "`"My phone number is `$.`""
$ws.Range("B36").Value = "6133992081, 2137584930, 2637490563, 2635478152, 613 399 2081, 213 758 4930, 263 749 0563, 263 547 8152, 613-399-2081, 213-758-4930, 263-749-0563, 263-547-8152"
$ws.Range("C36").Value = "PHONE_NUMBER"

$ws.Range("A37").Value = "`"My number is `$.`""
$ws.Range("B37").Value = "6133992081, 2137584930, 2637490563, 2635478152, 613 399 2081, 213 758 4930, 263 749 0563, 263 547 8152, 613-399-2081, 213-758-4930, 263-749-0563, 263-547-8152"
$ws.Range("C37").Value = "PHONE_NUMBER"

$ws.Range("A38").Value = "`"`$`""
$ws.Range("B38").Value = "6133992081, 2137584930, 2637490563, 2635478152, 613 399 2081, 213 758 4930, 263 749 0563, 263 547 8152, 613-399-2081, 213-758-4930, 263-749-0563, 263-547-8152"
$ws.Range("C38").Value = "PHONE_NUMBER"

# 4. Carry column-A's established cell format (s=1) down through the newly used rows 35-41
$ws.Range("A2").Copy()
$ws.Range("A35:A41").PasteSpecial(-4122)

# B38 is a brand-new row, so give it the same established format as B36:B37
$ws.Range("B36").Copy()
$ws.Range("B38").PasteSpecial(-4122)

# 5. Apply a new font (black Aptos Narrow) across B39:K41, matching a block formatting pass
$ws.Range("B39:K41").Font.Color = 0

$excel.CutCopyMode = 0

# 6. Move the active selection down to where the determiner now leaves off
$ws.Range("C43").Select()
